$d = $word.ActiveDocument

# 1. Uppercase the title text "Enunciado del Alcance" -> "ENUNCIADO DEL ALCANCE"
[void]$d.Content.Find.Execute("Enunciado del Alcance", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ENUNCIADO DEL ALCANCE", 2)

# 2. Insert three blank paragraphs right after the title paragraph (cover page grows)
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "ENUNCIADO DEL ALCANCE") {
        $titlePara = $p
        break
    }
}
if ($titlePara -eq $null) {
    throw "Could not locate the 'ENUNCIADO DEL ALCANCE' title paragraph"
}

$blankParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:afterAutospacing="0"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$nextParaIndex = $titlePara.Index + 1
for ($i = 0; $i -lt 3; $i++) {
    $nextPara = $d.Paragraphs($nextParaIndex)
    $insertPos = $nextPara.Range.Start
    $rng = $d.Range($insertPos, $insertPos)
    [void]$rng.InsertXML($blankParaXml)
    $nextParaIndex = $nextParaIndex + 1
}

# 3. Give the section a distinct title page (adds <w:titlePg/> to sectPr)
$d.Sections(1).PageSetup.DifferentFirstPageHeaderFooter = $true
